$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. These cells are stored as text
# (inline strings) in the workbook, not numbers.
$updates = @{
    "D2"  = "271.69"
    "D3"  = "22.65"
    "D4"  = "6.335"
    "D5"  = "0.06203"
    "D6"  = "3.656"
    "D7"  = "6.664"
    "D8"  = "1.401"
    "E8"  = "7FTXTokenFTTWorstin24h"
    "D9"  = "0.8322"
    "D10" = "0.01376"
    "D11" = "0.1612"
    "D12" = "0.08307"
    "D13" = "0.03555"
    "D14" = "0.03213"
    "D15" = "4.080"
    "D16" = "0.09303"
    "D17" = "0.001636"
    "D19" = "0.006353"
    "D20" = "0.005686"
    "D23" = "3.728"
    "D25" = "0.3338"
    "D27" = "0.0002707"
    "D40" = "0.04737"
    "D41" = "0.006943"
    "D42" = "0.003804"
    "E42" = "41CEJICEJI"
    "D44" = "0.01208"
    "D45" = "0.00006234"
    "D46" = "0.0009912"
    "D48" = "0.7830"
    "D49" = "0.002327"
    "D50" = "0.00002403"
}

foreach ($addr in $updates.Keys) {
    # Prefix with a leading apostrophe so Excel stores the value as literal
    # text (matching the workbook's existing inline-string cells) instead of
    # re-interpreting numeric-looking strings as numbers, which would strip
    # meaningful formatting such as trailing zeros (e.g. "0.06210", "1.400").
    $ws.Range($addr).Value = "'" + $updates[$addr]
}
